# "arreglamos el tablero final" - fix the final Jornada 6 / Jornada 7 columns
# on the "Jugadores" sheet: clean up the "J6 -  Minutos" / "J7 -  Minutos"
# header typos (double space) and fill in the Jornada 7 "Minutos" data that
# was missing, rolling the new minutes into the season "Minutos" totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jugadores")

# --- Fix the double-space typo in the J6 / J7 "Minutos" column headers ---
$ws.Range("AC1").Value = "J6 - Minutos"
$ws.Range("AE1").Value = "J7 - Minutos"

# --- Fill in the Jornada 7 minutes played (column AC) per player, and
#     update the season-total "Minutos" column (K) to include them. ---

# Hernán Barcos
$ws.Range("AC2").Value = 45
$ws.Range("K2").Value = 469

# Juan Freytes
$ws.Range("AC3").Value = 90
$ws.Range("K3").Value = 540

# Franco Saravia
$ws.Range("AC4").Value = 25
$ws.Range("K4").Value = 384

# Kevin Serna
$ws.Range("AC5").Value = 90
$ws.Range("K5").Value = 499

# Cecilio Waterman
$ws.Range("AC6").Value = 85
$ws.Range("K6").Value = 506

# Jiovany Ramos
$ws.Range("AC7").Value = 90
$ws.Range("K7").Value = 426

# Sebastián Rodríguez
$ws.Range("AC8").Value = 90
$ws.Range("K8").Value = 482

# Catriel Cabellos (no minutes played in J7)
$ws.Range("AC9").Value = 0

# Aldair Fuentes (no minutes played in J7)
$ws.Range("AC10").Value = 0

# Adrián Arregui
$ws.Range("AC12").Value = 90
$ws.Range("K12").Value = 418

# Jhamir D'Arrigo
$ws.Range("AC13").Value = 69
$ws.Range("K13").Value = 186

# Jesús Castillo
$ws.Range("AC14").Value = 45
$ws.Range("K14").Value = 158

# Renzo Garces
$ws.Range("AC15").Value = 90
$ws.Range("K15").Value = 342

# Gabriel Costa
$ws.Range("AC16").Value = 30
$ws.Range("K16").Value = 232

# Franco Zanelatto
$ws.Range("AC17").Value = 21
$ws.Range("K17").Value = 106

# Axel Moyano (no minutes played in J7)
$ws.Range("AC18").Value = 0

# Marco Huaman (no minutes played in J7)
$ws.Range("AC19").Value = 0

# Ángel De la Cruz
$ws.Range("AC21").Value = 60
$ws.Range("K21").Value = 150

# Sebastian Aranda (no minutes played in J7)
$ws.Range("AC27").Value = 0

# Víctor Guzmán
$ws.Range("AC28").Value = 9
$ws.Range("K28").Value = 19

# --- Leave the cursor where the author left it when done editing ---
$ws.Activate()
$ws.Range("R19").Select()
